# Duplicate "Sheet2" twice, inserting each copy directly after Sheet2 (i.e.
# before Sheet3), mirroring Excel's "Move or Copy... > Create a copy" flow.
# Excel auto-names the copies "Sheet2 (2)" and "Sheet2 (3)".
$wb = $excel.ActiveWorkbook
$sheet2 = $wb.Worksheets.Item("Sheet2")

$sheet2.Copy($null, $sheet2)
$sheet2.Copy($null, $wb.Worksheets.Item("Sheet2 (2)"))

# Land on the last copy ("Sheet2 (3)") as the active sheet/tab.
$wb.Worksheets.Item("Sheet2 (3)").Activate()
